$wb = $excel.ActiveWorkbook

# Switch to the second sheet ("This_is_sheet2") and make it the active sheet
$ws = $wb.Worksheets.Item("This_is_sheet2")
$ws.Activate()

# Populate the header row
$ws.Range("A1").Value = "First name"
$ws.Range("B1").Value = "Last name"
